$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list data: row number -> column letter -> new cell text
$updates = @{
    2 = @{ D = "43.783.73"; E = "  -0.22%  " }
    3 = @{ D = "2.312.40"; E = "  +3.79%  " }
    4 = @{ E = "  +0.16%  " }
    5 = @{ D = "269.16"; E = "  -0.27%  " }
    6 = @{ D = "92.61"; E = "  +6.48%  " }
    7 = @{ D = "0.628"; E = "  +1.84%  " }
    8 = @{ E = "  +0.07%  " }
    9 = @{ D = "0.619"; E = "  +1.84%  " }
    10 = @{ D = "44.83"; E = "  -3.22%  " }
    11 = @{ D = "0.0934"; E = "  +1.05%  " }
    12 = @{ D = "7.98"; E = "  +5.73%  " }
    13 = @{ E = "  +0.17%  " }
    14 = @{ D = "2.661.48"; E = "  +4.00%  " }
    15 = @{ D = "15.27"; E = "  +3.69%  " }
    16 = @{ D = "0.851"; E = "  +7.61%  " }
    17 = @{ D = "2.329.35"; E = "  +4.66%  " }
    18 = @{ D = "43.771.75"; E = "  -0.14%  " }
    19 = @{ E = "  +1.93%  " }
    20 = @{ D = "6.27"; E = "  +3.78%  " }
    21 = @{ D = "71.16"; E = "  +1.66%  " }
    22 = @{ D = "240.56"; E = "  +3.34%  " }
    23 = @{ E = "  -4.76%  " }
    24 = @{ D = "9.70"; E = "  +8.90%  " }
    25 = @{ E = "  -0.09%  " }
    26 = @{ E = "  -7.97%  " }
    27 = @{ E = "  +3.83%  " }
    28 = @{ D = "2.32"; E = "  +3.23%  " }
    29 = @{ D = "3.38"; E = "  -4.27%  " }
    30 = @{ D = "38.95"; E = "  -1.43%  " }
    31 = @{ D = "22.51"; E = "  +9.17%  " }
    32 = @{ D = "171.93"; E = "  -2.00%  " }
    33 = @{ D = "0.0890"; E = "  -1.29%  " }
    34 = @{ D = "5.53"; E = "  +1.69%  " }
    35 = @{ E = "  +1.37%  " }
    36 = @{ D = "0.110"; E = "  -0.76%  " }
    37 = @{ D = "4.53"; E = "  +3.50%  " }
    38 = @{ D = "0.0347"; E = "  -3.15%  " }
    39 = @{ D = "3.36"; E = "  +0.30%  " }
    40 = @{ E = "  +15.96%  " }
    41 = @{ D = "2.29"; E = "  +8.29%  " }
    42 = @{ D = "12.16"; E = "  -2.50%  " }
    43 = @{ D = "1.31"; E = "  +16.02%  " }
    44 = @{ D = "5.42"; E = "  +1.42%  " }
    45 = @{ D = "61.15"; E = "  -6.10%  " }
    46 = @{ D = "8.90"; E = "  +6.81%  " }
    47 = @{ E = "  +2.53%  " }
    48 = @{ D = "100.17"; E = "  -0.81%  " }
    49 = @{ D = "1.19"; E = "  -1.68%  " }
    50 = @{ D = "2.539.17"; E = "  +4.17%  " }
    51 = @{ D = "0.429"; E = "  -3.05%  " }
}

foreach ($row in $updates.Keys) {
    $entry = $updates[$row]
    foreach ($col in $entry.Keys) {
        $cell = $ws.Range("$col$row")
        # Force the cell to remain plain text so values such as "269.16"
        # or "9.70" are not silently coerced into numbers by Excel,
        # then restore the default (unstyled) cell style/format.
        $cell.NumberFormat = "@"
        $cell.Value = $entry[$col]
        $cell.Style = "Normal"
    }
}
